$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the changed range so numeric-looking
# strings (e.g. "212.04", "1.00") are stored as text, matching the
# original inline-string cell contents instead of being coerced to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.821.34'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '1.604.94'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '212.04'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('D12').Value = '1.830.21'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '1.610.35'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').Value = '0.528'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').Value = '65.28'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').Value = '0.0₃0744'
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').Value = '210.19'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '1.00'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '7.17'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = '4.32'
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('D22').Value = '2.26'
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('D23').Value = '9.06'
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('D24').Value = '144.08'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('D25').Value = '1.01'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').Value = '7.17'
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').Value = '15.39'
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('D29').Value = '0.0510'
$ws.Range('E29').Value = '  -1.83%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').Value = '3.28'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').Value = '1.28'
$ws.Range('E33').Value = '  +19.03%  '
$ws.Range('D34').Value = '1.282.29'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').Value = '2.48'
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').Value = '0.596'
$ws.Range('E37').Value = '  -3.47%  '
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('D39').Value = '0.829'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '5.48'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '2.18'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('D43').Value = '62.97'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').Value = '1.741.70'
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('D45').Value = '90.63'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').Value = '1.58'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('E47').Value = '  +2.47%  '
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.54'
$ws.Range('E49').Value = '  +2.85%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.399'
$ws.Range('E51').Value = '  +1.52%  '

# Restore the default (unstyled) appearance so only cell contents changed.
$ws.Range("B2:E51").Style = "Normal"
